$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("F2").Value = 1.53
$ws.Range("G2").Value = 1.6
$ws.Range("H2").Value = 5.9
$ws.Range("I2").Value = 6.8
$ws.Range("J2").Value = 4.6
$ws.Range("K2").Value = 5.2
$ws.Range("N2").Value = 4.1
$ws.Range("O2").Value = 1.26
$ws.Range("P2").Value = 2.1
$ws.Range("Q2").Value = 1.79
$ws.Range("T2").Value = 1.87
$ws.Range("U2").Value = 1.97
$ws.Range("W2").Value = 2.6
$ws.Range("X2").Value = 32
$ws.Range("Y2").Value = 1000
$ws.Range("AC2").Value = 13
$ws.Range("AD2").Value = 1000
$ws.Range("AF2").Value = 19
$ws.Range("AK2").Value = 50
$ws.Range("AN2").Value = 8.6
$ws.Range("G3").Value = 8.6
$ws.Range("H3").Value = 1.41
$ws.Range("J3").Value = 5.2
$ws.Range("K3").Value = 6.2
$ws.Range("L3").Value = 1.3
$ws.Range("T3").Value = 1.86
$ws.Range("U3").Value = 2.04
$ws.Range("AA3").Value = 13
$ws.Range("AB3").Value = 32
$ws.Range("AG3").Value = 85
$ws.Range("AI3").Value = 34
$ws.Range("F4").Value = 1.3
$ws.Range("H4").Value = 4.1
$ws.Range("I4").Value = 34
$ws.Range("J4").Value = 3.75
$ws.Range("K4").Value = 26
$ws.Range("L4").Value = 1.01
$ws.Range("R4").Value = 1.19
$ws.Range("U4").Value = 2.02
$ws.Range("W4").Value = 2.66
$ws.Range("AC4").Value = 100
$ws.Range("AJ4").Value = 190
$ws.Range("F5").Value = 2.16
$ws.Range("G5").Value = 2.38
$ws.Range("K5").Value = 3.65
$ws.Range("L5").Value = 1.45
$ws.Range("N5").Value = 3.15
$ws.Range("O5").Value = 1.37
$ws.Range("R5").Value = 1.28
$ws.Range("S5").Value = 3.85
$ws.Range("T5").Value = 1.84
$ws.Range("Y5").Value = 14
$ws.Range("AA5").Value = 95
$ws.Range("AB5").Value = 10.5
$ws.Range("AC5").Value = 9
$ws.Range("AD5").Value = 18
$ws.Range("AE5").Value = 55
$ws.Range("AF5").Value = 16.5
$ws.Range("AI5").Value = 190
$ws.Range("AN5").Value = 26
$ws.Range("G6").Value = 1.87
$ws.Range("H6").Value = 5.4
$ws.Range("J6").Value = 3.4
$ws.Range("L6").Value = 1.39
$ws.Range("N6").Value = 3.1
$ws.Range("O6").Value = 1.39
$ws.Range("P6").Value = 1.72
$ws.Range("Q6").Value = 2.12
$ws.Range("R6").Value = 1.26
$ws.Range("S6").Value = 4
$ws.Range("T6").Value = 1.97
$ws.Range("U6").Value = 1.86
$ws.Range("W6").Value = 2.14
$ws.Range("X6").Value = 12.5
$ws.Range("Y6").Value = 990
$ws.Range("Z6").Value = 1000
$ws.Range("AB6").Value = 7.4
$ws.Range("AC6").Value = 8.4
$ws.Range("AD6").Value = 60
$ws.Range("AF6").Value = 11
$ws.Range("AG6").Value = 9
$ws.Range("AK6").Value = 65
$ws.Range("AL6").Value = 1000
$ws.Range("AN6").Value = 55
